$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Max value for ca_distances_calc (C2): 7.1 -> 6.9
$ws.Range("C2").Value = 6.9

# Update the Min value for ratio (B5): 1 -> 0.95
$ws.Range("B5").Value = 0.95

# Move/restore the active selection to C11 (matches saved cursor position)
$ws.Range("C11").Select() | Out-Null
